# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value2 = "Datos actualizados a 25 de Abril de 2020 a las 21:22"

# Swap India/Peru ordering: India overtakes Peru in total cases, so India's
# updated stats now occupy row 19 and Peru's (unchanged) stats move to row 20.
# Row 19 country name stays as shared-string index 23, but that string's text
# changes from "Peru" to "India"; row 20 keeps index 24, changing from "India"
# to "Peru". We just set the cell text directly and Excel will manage the
# shared string table.
$ws.Range("A19").Value2 = "India"
$ws.Range("A20").Value2 = "Peru"

# Row 4 - Estados Unidos
$ws.Range("B4").Value2 = 946921
$ws.Range("C4").Value2 = 21689
$ws.Range("D4").Value2 = 115910
$ws.Range("E4").Value2 = 777550
$ws.Range("F4").Value2 = 15100
$ws.Range("G4").Value2 = 1268
$ws.Range("H4").Value2 = 53461

# Row 14 - Brasil
$ws.Range("B14").Value2 = 57382
$ws.Range("C14").Value2 = 4387
$ws.Range("D14").Value2 = 27655
$ws.Range("E14").Value2 = 25809
$ws.Range("F14").Value2 = 8318
$ws.Range("G14").Value2 = 248
$ws.Range("H14").Value2 = 3918

# Row 19 - now India (updated stats)
$ws.Range("B19").Value2 = 26283
$ws.Range("C19").Value2 = 1836
$ws.Range("D19").Value2 = 5939
$ws.Range("E19").Value2 = 19519
$ws.Range("F19").Value2 = 0
$ws.Range("G19").Value2 = 45
$ws.Range("H19").Value2 = 825

# Row 20 - now Peru (stats carried over unchanged from old row 19)
$ws.Range("B20").Value2 = 25331
$ws.Range("C20").Value2 = 3683
$ws.Range("D20").Value2 = 7797
$ws.Range("E20").Value2 = 16834
$ws.Range("F20").Value2 = 545
$ws.Range("G20").Value2 = 66
$ws.Range("H20").Value2 = 700

# Row 98 - Costa Rica
$ws.Range("B98").Value2 = 693
$ws.Range("C98").Value2 = 6
$ws.Range("D98").Value2 = 242
$ws.Range("E98").Value2 = 445
$ws.Range("F98").Value2 = 7
$ws.Range("G98").Value2 = 0
$ws.Range("H98").Value2 = 6

# Row 106 - Estado de Palestina
$ws.Range("B106").Value2 = 495
$ws.Range("C106").Value2 = 11
$ws.Range("D106").Value2 = 92
$ws.Range("E106").Value2 = 399
$ws.Range("F106").Value2 = 0
$ws.Range("G106").Value2 = 0
$ws.Range("H106").Value2 = 4
